$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy style from B1 to C1 first so the new header cell matches formatting
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Update header row text
$ws.Range("B1").Value = "cases1"
$ws.Range("C1").Value = "cases2"

# Copy column B values (rows 2-64) into column C
for ($r = 2; $r -le 64; $r++) {
    $ws.Cells.Item($r, 3).Value = $ws.Cells.Item($r, 2).Value2
}

$ws.Range("C2").Select()
